$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26; this pushes the existing rows 26-44 down to 27-45.
$ws.Rows("26").Insert()

# The row that used to be row 26 is now row 27. Duplicate it back into the
# freshly-inserted row 26 so we have a full copy of that record there too.
$srcVals = $ws.Range("A27:T27").Value()
$ws.Range("A26:T26").Value = $srcVals

# Now adjust the two fields that differ for this new record: the date (column D)
# and the volume (column M).
$ws.Range("D26").Value = 44839
$ws.Range("M26").Value = 15
